# Update the "Waypart" (column G) raw sample counts for each benchmark.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1597567
$ws.Range("G4").Value = 20296155
$ws.Range("G5").Value = 20548493
$ws.Range("G6").Value = 1753649

# Flip the normalization formulas in rows 11-12 so they divide the
# baseline by each column instead of each column by the baseline.
$ws.Range("F11").Formula = "=`$F`$8/F8"
$ws.Range("G11").Formula = "=`$F`$8/G8"
$ws.Range("H11").Formula = "=`$F`$8/H8"
$ws.Range("I11").Formula = "=`$F`$8/I8"

$ws.Range("F12").Formula = "=`$F`$9/F9"
$ws.Range("G12").Formula = "=`$F`$9/G9"
$ws.Range("H12").Formula = "=`$F`$9/H9"
$ws.Range("I12").Formula = "=`$F`$9/I9"

# Move the active selection to match the author's final cursor position.
$ws.Range("I14").Select()
